# Auto-generated edit script: update res_bus vm_pu values for 380kV case (B2:N25, excluding G/H)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.021017892867217
$ws.Range("D2").Value = 1.030200767296413
$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("F2").Value = 1.036876815792654
$ws.Range("I2").Value = 1.030225492146027
$ws.Range("J2").Value = 1.026211730096811
$ws.Range("K2").Value = 1.033012466510578
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("M2").Value = 1.039669303220529
$ws.Range("N2").Value = 1.012829858124939
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.021837642343808
$ws.Range("D3").Value = 1.030826974434112
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("F3").Value = 1.037688220676287
$ws.Range("I3").Value = 1.030346042488668
$ws.Range("J3").Value = 1.026669593382548
$ws.Range("K3").Value = 1.033447468386647
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("M3").Value = 1.040290395545401
$ws.Range("N3").Value = 1.012980575344422
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.022368297281582
$ws.Range("D4").Value = 1.031231953556889
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("F4").Value = 1.038213581891308
$ws.Range("I4").Value = 1.030422151984871
$ws.Range("J4").Value = 1.026965440425052
$ws.Range("K4").Value = 1.033728036259707
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("M4").Value = 1.040691961526404
$ws.Range("N4").Value = 1.013077953889465
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.02259143623453
$ws.Range("D5").Value = 1.031402152508963
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("F5").Value = 1.038434520280707
$ws.Range("I5").Value = 1.030453693884817
$ws.Range("J5").Value = 1.027089712437944
$ws.Range("K5").Value = 1.033845768282142
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("M5").Value = 1.040860700989619
$ws.Range("N5").Value = 1.013118856549014
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.022628905227415
$ws.Range("D6").Value = 1.031430726424073
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("F6").Value = 1.038471621238505
$ws.Range("I6").Value = 1.030458963226799
$ws.Range("J6").Value = 1.027110572232544
$ws.Range("K6").Value = 1.033865523138476
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("M6").Value = 1.040889028412546
$ws.Range("N6").Value = 1.013125722204108
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.022371278674363
$ws.Range("D7").Value = 1.031234227979271
$ws.Range("E7").Value = 0.9943035907982488
$ws.Range("F7").Value = 1.038216533782668
$ws.Range("I7").Value = 1.030422575236491
$ws.Range("J7").Value = 1.026967101356492
$ws.Range("K7").Value = 1.033729610262004
$ws.Range("L7").Value = 0.9968970624462087
$ws.Range("M7").Value = 1.040694216543473
$ws.Range("N7").Value = 1.013078500571977
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.021294884081322
$ws.Range("D8").Value = 1.0304124410901
$ws.Range("E8").Value = 0.9929600610674301
$ws.Range("F8").Value = 1.037150964764495
$ws.Range("I8").Value = 1.030266624425034
$ws.Range("J8").Value = 1.026366553573644
$ws.Range("K8").Value = 1.033159664438271
$ws.Range("L8").Value = 0.995817528259106
$ws.Range("M8").Value = 1.039879269979227
$ws.Range("N8").Value = 1.012880823599049
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.019399910744079
$ws.Range("D9").Value = 1.028962755493903
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.035275895585319
$ws.Range("I9").Value = 1.029977350269543
$ws.Range("J9").Value = 1.025305142060444
$ws.Range("K9").Value = 1.032148461865461
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.038440837295449
$ws.Range("N9").Value = 1.012531395495599
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.018137879137439
$ws.Range("D10").Value = 1.027995339849007
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.034027714547674
$ws.Range("I10").Value = 1.02977482852293
$ws.Range("J10").Value = 1.024595477843943
$ws.Range("K10").Value = 1.031469784055265
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.037480374697414
$ws.Range("N10").Value = 1.012297731829471
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.01759172937188
$ws.Range("D11").Value = 1.027576234254459
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.033487705283266
$ws.Range("I11").Value = 1.029684851812605
$ws.Range("J11").Value = 1.024287713161909
$ws.Range("K11").Value = 1.031174849628369
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.037064147388405
$ws.Range("N11").Value = 1.012196389382907
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.017388914370743
$ws.Range("D12").Value = 1.027420530362999
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.033287193074027
$ws.Range("I12").Value = 1.029651088290091
$ws.Range("J12").Value = 1.024173325464436
$ws.Range("K12").Value = 1.031065139654946
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.036909492534415
$ws.Range("N12").Value = 1.012158721984943
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.017432416619793
$ws.Range("D13").Value = 1.027453930677309
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.033330200366794
$ws.Range("I13").Value = 1.02965834614967
$ws.Range("J13").Value = 1.02419786517044
$ws.Range("K13").Value = 1.031088679952511
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.036942668738583
$ws.Range("N13").Value = 1.01216680286264
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.0175749636053
$ws.Range("D14").Value = 1.027563364311443
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.033471129416189
$ws.Range("I14").Value = 1.029682067884216
$ws.Range("J14").Value = 1.024278259263935
$ws.Range("K14").Value = 1.031165784186655
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.03705136456999
$ws.Range("N14").Value = 1.01219327627736
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.017662798091217
$ws.Range("D15").Value = 1.027630786143523
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.03355796995798
$ws.Range("I15").Value = 1.029696638314207
$ws.Range("J15").Value = 1.024327783453973
$ws.Range("K15").Value = 1.031213269730197
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.037118329150236
$ws.Range("N15").Value = 1.012209584214272
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.018174132134124
$ws.Range("D16").Value = 1.028023150279762
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.034063563072016
$ws.Range("I16").Value = 1.029780751947789
$ws.Range("J16").Value = 1.024615893261515
$ws.Range("K16").Value = 1.031489335619065
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.037507991299438
$ws.Range("N16").Value = 1.01230445416307
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.018494964752313
$ws.Range("D17").Value = 1.02826921551035
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.034380833181456
$ws.Range("I17").Value = 1.029832903497001
$ws.Range("J17").Value = 1.024796490526627
$ws.Range("K17").Value = 1.031662221017557
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.037752326055044
$ws.Range("N17").Value = 1.012363919820521
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.018682131594651
$ws.Range("D18").Value = 1.028412720979479
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.034565935852155
$ws.Range("I18").Value = 1.029863102232916
$ws.Range("J18").Value = 1.024901783905106
$ws.Range("K18").Value = 1.031762959538628
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.037894809402467
$ws.Range("N18").Value = 1.012398589229004
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.018745955795594
$ws.Range("D19").Value = 1.028461649149715
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.034629058547834
$ws.Range("I19").Value = 1.029873361813721
$ws.Range("J19").Value = 1.024937678394552
$ws.Range("K19").Value = 1.031797291323564
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.037943386896272
$ws.Range("N19").Value = 1.01241040789109
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.018460539272124
$ws.Range("D20").Value = 1.028242817117413
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.034346788496328
$ws.Range("I20").Value = 1.0298273309178
$ws.Range("J20").Value = 1.024777118893228
$ws.Range("K20").Value = 1.031643682653247
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.037726114647402
$ws.Range("N20").Value = 1.012357541360188
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.017532985696719
$ws.Range("D21").Value = 1.027531139630708
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.033429627340445
$ws.Range("I21").Value = 1.02967509186292
$ws.Range("J21").Value = 1.024254587134935
$ws.Range("K21").Value = 1.031143083267018
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.037019357709753
$ws.Range("N21").Value = 1.012185481184031
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.016950082163425
$ws.Range("D22").Value = 1.027083510493651
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.03285338506269
$ws.Range("I22").Value = 1.029577393899408
$ws.Range("J22").Value = 1.023925645179295
$ws.Range("K22").Value = 1.030827421995453
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.036574706013775
$ws.Range("N22").Value = 1.012077159729864
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.017259062716822
$ws.Range("D23").Value = 1.027320822589485
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.033158822164106
$ws.Range("I23").Value = 1.029629372748085
$ws.Range("J23").Value = 1.024100061560833
$ws.Range("K23").Value = 1.030994846224218
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.036810450798803
$ws.Range("N23").Value = 1.012134596150458
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.018476094570386
$ws.Range("D24").Value = 1.028254745477556
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.034362171687183
$ws.Range("I24").Value = 1.029829849607499
$ws.Range("J24").Value = 1.024785872242631
$ws.Range("K24").Value = 1.031652059659294
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.03773795855635
$ws.Range("N24").Value = 1.012360423561039
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.019889587257796
$ws.Range("D25").Value = 1.029337711077012
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.035760325734746
$ws.Range("I25").Value = 1.030053843442111
$ws.Range("J25").Value = 1.025579910291318
$ws.Range("K25").Value = 1.032410688961113
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.038812979708511
$ws.Range("N25").Value = 1.012621858305038
